# Scheduled data refresh: update market-board derived profit figures
# (currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ /
#  LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ columns, i.e.
#  H:N) across the crafting-class worksheets with freshly pulled values.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 178
$ws.Range("I2").Value = 83.333336
$ws.Range("J2").Value = 234.8
$ws.Range("K2").Value = 83.333336
$ws.Range("L2").Value = 234.8
$ws.Range("M2").Value = 29.666664
$ws.Range("N2").Value = -460.8

$ws.Range("H4").Value = 293.77777
$ws.Range("I4").Value = 158.66667
$ws.Range("J4").Value = 564
$ws.Range("K4").Value = 158.66667
$ws.Range("L4").Value = 564
$ws.Range("M4").Value = -44.66667000000001
$ws.Range("N4").Value = -792

$ws.Range("H5").Value = 119.7
$ws.Range("I5").Value = 132.14285
$ws.Range("J5").Value = 90.666664
$ws.Range("K5").Value = 132.14285
$ws.Range("L5").Value = 90.666664
$ws.Range("M5").Value = -17.14285000000001
$ws.Range("N5").Value = -320.666664

$ws.Range("H8").Value = 490.73685
$ws.Range("I8").Value = 490.73685
$ws.Range("K8").Value = 1472.21055
$ws.Range("M8").Value = -1333.21055

$ws.Range("H9").Value = 123.92308
$ws.Range("I9").Value = 168.375
$ws.Range("J9").Value = 52.8
$ws.Range("K9").Value = 168.375
$ws.Range("L9").Value = 52.8
$ws.Range("M9").Value = 0.625
$ws.Range("N9").Value = -390.8

$ws.Range("H10").Value = 1500
$ws.Range("I10").Value = 1000
$ws.Range("J10").Value = 2000
$ws.Range("K10").Value = 1000
$ws.Range("L10").Value = 2000
$ws.Range("M10").Value = -707
$ws.Range("N10").Value = -2586

$ws.Range("H12").Value = 76923480
$ws.Range("I12").Value = 380
$ws.Range("J12").Value = 142857570
$ws.Range("K12").Value = 380
$ws.Range("L12").Value = 142857570
$ws.Range("M12").Value = -210
$ws.Range("N12").Value = -142857910

$ws.Range("H13").Value = 41570.668
$ws.Range("J13").Value = 41570.668
$ws.Range("L13").Value = 41570.668
$ws.Range("N13").Value = -41908.668

$ws.Range("H16").Value = 3115.2856
$ws.Range("I16").Value = 3081.75
$ws.Range("J16").Value = 3160
$ws.Range("K16").Value = 3081.75
$ws.Range("L16").Value = 3160
$ws.Range("M16").Value = -2851.75
$ws.Range("N16").Value = -3620

$ws.Range("H19").Value = 650.13513
$ws.Range("I19").Value = 573.2941
$ws.Range("J19").Value = 715.45
$ws.Range("K19").Value = 573.2941
$ws.Range("L19").Value = 715.45
$ws.Range("M19").Value = -398.2941
$ws.Range("N19").Value = -1065.45

$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()

$ws.Range("H21").Value = 24803.9
$ws.Range("I21").Value = 33008.75
$ws.Range("K21").Value = 33008.75
$ws.Range("M21").Value = -32540.75

$ws.Range("H23").Value = 24803.9
$ws.Range("I23").Value = 33008.75
$ws.Range("K23").Value = 33008.75
$ws.Range("M23").Value = -32774.75

$ws.Range("H29").Value = 885.6667
$ws.Range("I29").Value = 645
$ws.Range("J29").Value = 2450
$ws.Range("K29").Value = 1935
$ws.Range("L29").Value = 7350
$ws.Range("M29").Value = -1654
$ws.Range("N29").Value = -7912

$ws.Range("H31").Value = 821.8333
$ws.Range("I31").Value = 386
$ws.Range("J31").Value = 3001
$ws.Range("K31").Value = 1158
$ws.Range("L31").Value = 9003
$ws.Range("M31").Value = -928
$ws.Range("N31").Value = -9463

$ws.Range("H34").Value = 14158.667
$ws.Range("I34").Value = 11970
$ws.Range("J34").Value = 14596.4
$ws.Range("K34").Value = 11970
$ws.Range("L34").Value = 14596.4
$ws.Range("M34").Value = -11767
$ws.Range("N34").Value = -15002.4

$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").ClearContents()

$ws.Range("H36").Value = 14158.667
$ws.Range("I36").Value = 11970
$ws.Range("J36").Value = 14596.4
$ws.Range("K36").Value = 11970
$ws.Range("L36").Value = 14596.4
$ws.Range("M36").Value = -11255
$ws.Range("N36").Value = -16026.4

$ws.Range("H38").Value = 349.23077
$ws.Range("I38").Value = 130
$ws.Range("J38").Value = 700
$ws.Range("K38").Value = 390
$ws.Range("L38").Value = 2100
$ws.Range("M38").Value = -18
$ws.Range("N38").Value = -2844

$ws.Range("H39").Value = 192.06667
$ws.Range("I39").Value = 31.6
$ws.Range("J39").Value = 513
$ws.Range("K39").Value = 94.80000000000001
$ws.Range("L39").Value = 1539
$ws.Range("M39").Value = 201.2
$ws.Range("N39").Value = -2131

$ws.Range("H40").Value = 1562.125
$ws.Range("I40").Value = 1111.5454
$ws.Range("J40").Value = 1943.3846
$ws.Range("K40").Value = 1111.5454
$ws.Range("L40").Value = 1943.3846
$ws.Range("M40").Value = -936.5454
$ws.Range("N40").Value = -2293.3846

$ws.Range("H41").Value = 968.125
$ws.Range("I41").Value = 1107.5
$ws.Range("J41").Value = 550
$ws.Range("K41").Value = 1107.5
$ws.Range("L41").Value = 550
$ws.Range("M41").Value = -667.5
$ws.Range("N41").Value = -1430

$ws.Range("H42").Value = 95
$ws.Range("I42").Value = 32.5
$ws.Range("J42").Value = 157.5
$ws.Range("K42").Value = 97.5
$ws.Range("L42").Value = 472.5
$ws.Range("M42").Value = 132.5
$ws.Range("N42").Value = -932.5

$ws.Range("H43").Value = 3251
$ws.Range("I43").Value = 2500
$ws.Range("J43").Value = 4002
$ws.Range("K43").Value = 2500
$ws.Range("L43").Value = 4002
$ws.Range("M43").Value = -2431
$ws.Range("N43").Value = -4140

$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()

$ws.Range("H48").Value = 5000
$ws.Range("J48").Value = 5000
$ws.Range("L48").Value = 15000
$ws.Range("N48").Value = -15584

$ws.Range("H56").Value = 5000
$ws.Range("J56").Value = 5000
$ws.Range("L56").Value = 15000
$ws.Range("N56").Value = -16068

$ws.Range("H100").Value = 2156.9285
$ws.Range("I100").Value = 1657.909
$ws.Range("J100").Value = 3986.6667
$ws.Range("K100").Value = 1657.909
$ws.Range("L100").Value = 3986.6667
$ws.Range("M100").Value = -1116.909
$ws.Range("N100").Value = -5068.6667

$ws.Range("H132").Value = 2832.6416
$ws.Range("I132").Value = 2425.0732
$ws.Range("J132").Value = 4225.1665
$ws.Range("K132").Value = 7275.219599999999
$ws.Range("L132").Value = 12675.4995
$ws.Range("M132").Value = -4745.219599999999
$ws.Range("N132").Value = -17735.4995

$ws.Range("H135").Value = 658.25
$ws.Range("I135").Value = 437.33334
$ws.Range("J135").Value = 879.1667
$ws.Range("K135").Value = 3936.00006
$ws.Range("L135").Value = 7912.5003
$ws.Range("M135").Value = -1401.00006
$ws.Range("N135").Value = -12982.5003

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5095.45
$ws.Range("I32").Value = 4683.652
$ws.Range("J32").Value = 8427.272000000001
$ws.Range("K32").Value = 4683.652
$ws.Range("L32").Value = 8427.272000000001
$ws.Range("M32").Value = -4396.652
$ws.Range("N32").Value = -9001.272000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 128991.875
$ws.Range("I20").Value = 1004
$ws.Range("J20").Value = 171654.5
$ws.Range("K20").Value = 1004
$ws.Range("L20").Value = 171654.5
$ws.Range("M20").Value = -757
$ws.Range("N20").Value = -172148.5

$ws.Range("H120").Value = 22253.666
$ws.Range("J120").Value = 22253.666
$ws.Range("L120").Value = 22253.666
$ws.Range("N120").Value = -31929.666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2100.61
$ws.Range("I31").Value = 1180.5962
$ws.Range("J31").Value = 3097.2917
$ws.Range("K31").Value = 1180.5962
$ws.Range("L31").Value = 3097.2917
$ws.Range("M31").Value = -885.5962
$ws.Range("N31").Value = -3687.2917

$ws.Range("H34").Value = 2100.61
$ws.Range("I34").Value = 1180.5962
$ws.Range("J34").Value = 3097.2917
$ws.Range("K34").Value = 1180.5962
$ws.Range("L34").Value = 3097.2917
$ws.Range("M34").Value = -978.5962
$ws.Range("N34").Value = -3501.2917

$ws.Range("H94").Value = 1978.8
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 1978.8
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 1978.8
$ws.Range("M94").ClearContents()
$ws.Range("N94").Value = -2880.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 735.7568
$ws.Range("I5").Value = 541.8
$ws.Range("J5").Value = 963.94116
$ws.Range("K5").Value = 1625.4
$ws.Range("L5").Value = 2891.82348
$ws.Range("M5").Value = -1513.4
$ws.Range("N5").Value = -3115.82348

$ws.Range("H122").Value = 947.02856
$ws.Range("I122").Value = 478.375
$ws.Range("K122").Value = 4305.375
$ws.Range("M122").Value = -1855.375

$ws.Range("H134").Value = 55841.047
$ws.Range("I134").Value = 67774.234
$ws.Range("J134").Value = 5125
$ws.Range("K134").Value = 203322.702
$ws.Range("L134").Value = 15375
$ws.Range("M134").Value = -198252.702
$ws.Range("N134").Value = -25515

$ws.Range("H135").Value = 735.7568
$ws.Range("I135").Value = 541.8
$ws.Range("J135").Value = 963.94116
$ws.Range("K135").Value = 4876.2
$ws.Range("L135").Value = 8675.470439999999
$ws.Range("M135").Value = -2341.2
$ws.Range("N135").Value = -13745.47044

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 47.705883
$ws.Range("I2").Value = 38.42857
$ws.Range("J2").Value = 54.2
$ws.Range("K2").Value = 38.42857
$ws.Range("L2").Value = 54.2
$ws.Range("M2").Value = 74.57142999999999
$ws.Range("N2").Value = -280.2

$ws.Range("H113").Value = 556827.3
$ws.Range("I113").Value = 834411
$ws.Range("J113").Value = 1660
$ws.Range("K113").Value = 834411
$ws.Range("L113").Value = 1660
$ws.Range("M113").Value = -832241
$ws.Range("N113").Value = -6000

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 3148.4849
$ws.Range("I2").Value = 450
$ws.Range("J2").Value = 3322.5806
$ws.Range("K2").Value = 450
$ws.Range("L2").Value = 3322.5806
$ws.Range("M2").Value = -338
$ws.Range("N2").Value = -3546.5806

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 100790.3
$ws.Range("I113").Value = 200640.4
$ws.Range("J113").Value = 940.2
$ws.Range("K113").Value = 601921.2
$ws.Range("L113").Value = 2820.6
$ws.Range("M113").Value = -599751.2
$ws.Range("N113").Value = -7160.6

$ws.Range("H132").Value = 2512.4043
$ws.Range("I132").Value = 2831.6865
$ws.Range("J132").Value = 1720.1111
$ws.Range("K132").Value = 8495.059499999999
$ws.Range("L132").Value = 5160.3333
$ws.Range("M132").Value = -5965.059499999999
$ws.Range("N132").Value = -10220.3333
